# "Redeem points 79174445 30.0"
#
# The source sheet (redemptions) has 4 data rows (A2:C5). This edit:
#   1. Normalizes A5 (phone number) from text to a real number, matching
#      the numeric storage used by the other phone cells in the column.
#   2. Appends a new redemption row (row 6): same phone number (stored as
#      text, like the original A5 was), 30 points redeemed, and a new
#      timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) A5 was stored as text ("79174445"); make it a genuine number like A2:A4.
$ws.Range("A5").Value = 79174445

# 2) Append the new redemption row.
#    A6 must stay a text cell (leading apostrophe forces text without
#    touching the cell's number format/style), B6 is numeric, C6 is text.
$ws.Range("A6").Value = "'79174445"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = "2025-08-18T08:57:38"
